# Git Notes.docx edit
#
#   "Added additional text file to test git"
#
# Appends two new text paragraphs (each preceded by a blank paragraph) after
# the existing final paragraph ("Pointer to tip of current branch in repo...")
# and relocates the document's "_GoBack" bookmark so it again sits right
# after the very last run in the document (Word always keeps _GoBack at the
# spot of the most recent edit).

$d = $word.ActiveDocument

function Get-DocEndRange {
    # Collapsed Range sitting right at the end of the last paragraph's text
    # (i.e. immediately before that paragraph's mark).
    $n = $d.Paragraphs.Count
    $p = $d.Paragraphs.Item($n)
    $r = $p.Range
    $r.Collapse(0)   # wdCollapseEnd
    return $r
}

# Make sure we're anchored on the paragraph that currently ends the story.
$probe = $d.Content
if (-not $probe.Find.Execute("Points where we")) {
    throw "Expected anchor text ('Points where we...') was not found."
}

# --- 1) blank paragraph right after "Pointer to tip ... left off." ---------
$r = Get-DocEndRange
$r.InsertParagraphAfter()

# --- 2) "You always have to git commit -m" ---------------------------------
$r = Get-DocEndRange
$r.InsertParagraphAfter()
$r = Get-DocEndRange
$r.InsertAfter("You always have to git commit -m")

# --- 3) blank paragraph -----------------------------------------------------
$r = Get-DocEndRange
$r.InsertParagraphAfter()

# --- 4) "Git status shows you that you have nothing to commit. But what if
#         we want to make some changes?" ------------------------------------
$r = Get-DocEndRange
$r.InsertParagraphAfter()
$r = Get-DocEndRange
$r.InsertAfter("Git status shows you that you have nothing to commit. But what if we want to make some changes?")

# --- 5) move the "_GoBack" bookmark to the new end of the document ---------
# Word always keeps a single "_GoBack" bookmark marking the location of the
# last edit, collapsed at that point. Remove it from its old spot (end of
# the original last paragraph) and re-create it collapsed at the new end of
# story. Bookmarks.Add on a collapsed Range that sits exactly at a paragraph
# boundary gets mis-anchored in this host, so we briefly insert a marker
# character to give the Range interior context, anchor the bookmark next to
# it, then delete the marker again.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$tail = Get-DocEndRange
$tail.InsertAfter("ZZGOBACKMARKERZZ")
$bm = $tail.Duplicate
$bm.Collapse(1)        # wdCollapseStart -> right before the marker text
$d.Bookmarks.Add("_GoBack", $bm)
$tail.Delete()

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
